$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Locate the "A9.com | Analytics Data Platform" text robustly (avoids
# relying on hard-coded character offsets).
# ----------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("A9.com | Analytics Data Platform", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$matchStart = $findRng.Start
$matchEnd   = $findRng.End

$oldCompanyLen = 6                      # Len("A9.com")
$posR1R2 = $matchStart + $oldCompanyLen # boundary right after "A9.com"
$posR2R3 = $matchEnd                    # boundary right after " | Analytics Data Platform"

# ----------------------------------------------------------------------
# Step 1: isolate the "A9.com" run from the " | Analytics Data Platform"
# run so the upcoming text replace doesn't get re-flowed/merged into its
# neighbor (both runs share identical formatting). A short-lived bookmark
# at the boundary acts as a hard split point for the engine's run
# retokenizer.
# ----------------------------------------------------------------------
$placeholder1 = $d.Range($posR1R2, $posR1R2)
$placeholder1.InsertBefore("Y")
$bm1Range = $d.Range($posR1R2, $posR1R2 + 1)
$d.Bookmarks.Add("zzTempSplit1", $bm1Range)
$bm1 = $d.Bookmarks("zzTempSplit1")
$bm1.Range.Text = ""

# ----------------------------------------------------------------------
# Step 2: rename the company, "A9.com" -> "Amazon" (run stays isolated
# thanks to the zzTempSplit1 bookmark still sitting at the boundary).
# ----------------------------------------------------------------------
$companyRng = $d.Range($matchStart, $matchStart + $oldCompanyLen)
$companyRng.Text = "Amazon"

# ----------------------------------------------------------------------
# Step 3: insert the new " at A9.com" run between the
# " | Analytics Data Platform" run and the ", Palo Alto, CA" run using
# the same placeholder+bookmark technique, so it lands in its own run
# rather than merging into a neighbor.
# ----------------------------------------------------------------------
$placeholder2 = $d.Range($posR2R3, $posR2R3)
$placeholder2.InsertBefore("Z")
$bm2Range = $d.Range($posR2R3, $posR2R3 + 1)
$d.Bookmarks.Add("zzTempSplit2", $bm2Range)
$bm2 = $d.Bookmarks("zzTempSplit2")
$bm2.Range.Text = " at A9.com"

# ----------------------------------------------------------------------
# Step 4: clean up. Drop the first helper bookmark entirely, and turn the
# second one into the zero-length "_GoBack" bookmark sitting right after
# the newly inserted " at A9.com" text (mirroring Word's own behavior of
# relocating _GoBack to the most recent edit position).
# ----------------------------------------------------------------------
$d.Bookmarks("zzTempSplit1").Delete()

$bm2Final = $d.Bookmarks("zzTempSplit2")
$goBackPos = $bm2Final.End
$d.Bookmarks("zzTempSplit2").Delete()

$goBackRng = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRng)

Write-Output "done"
